$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("Data")
$notes = $wb.Worksheets.Item("Notes")

# Row 6 on the Data sheet had its monthly figures shifted into the wrong
# columns (JUN/JUL/AUG instead of MAR/APR/MAY, and NOV/DEC instead of AUG/SEP).
# Fix the error by moving the formulas three columns to the left and
# zeroing out the columns that previously held them.

$data.Range("H6").Formula = "=312575.3/325851"
$data.Range("I6").Formula = "=285472.2/325851"
$data.Range("J6").Formula = "=137489.5/325851"
$data.Range("K6").Value = 0
$data.Range("L6").Value = 0
$data.Range("M6").Formula = "=36/325851"
$data.Range("N6").Formula = "=86.6/325851"
$data.Range("O6").Value = 0
$data.Range("P6").Value = 0
$data.Range("Q6").Value = 0

# Switch the active/selected sheet from Notes back to Data, and update the
# last selected cell on each sheet to match.
$notes.Range("A4").Select() | Out-Null
$data.Activate() | Out-Null
$data.Range("N6").Select() | Out-Null
